# Generate Report for Handoff
# Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
# and the corresponding "Latest HO Xliff Generate Date" / "Latest Handoff
# Datetime" timestamps are refreshed to reflect the new handoff generation.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (E, F) + HO Xliff generate date (G)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-23 17:04:32"

# zh-cn sheet: Status (C) + Latest Handoff Datetime (H)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-23 17:04:28"

# de-de sheet: Status (C) ; its Latest Handoff Datetime (H2) shares the same
# underlying text as the Overview's G2 ("2016-08-23 17:03:38") so it moves
# to the same refreshed value ("2016-08-23 17:04:32") automatically.
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-23 17:04:32"

# Column widths shrink because the new "Ready for handoff" status text is
# shorter than the previous "Handed back: in sync with en-US" text, so the
# autofit width for the Status columns changes from 29.9777047293527 to
# 17.2159881591797 (in character-width units). The host's ColumnWidth
# setter quantizes to whole-pixel increments, so 16.33 is the calibrated
# input that lands on the closest representable stored width.
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
